# IST price update 2025-12-21 22:44
#
# The tracker keeps one column per price-check timestamp, newest first
# (column B). A new check was just run, so we insert a fresh column right
# after the "SKU Name" column, push every existing snapshot one column to
# the right, stamp the new column's header with the new check time, and
# carry forward the prices found (unchanged from the previous check).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; B..S (old snapshots) shift to C..T.
$ws.Columns("B").Insert()

# Insert resets the new column's width to the sheet default; restore it to
# match the other snapshot columns (all 21 / ColumnWidth ~20.17).
$ws.Columns("B").ColumnWidth = $ws.Columns("C").ColumnWidth

# Header for the new, most-recent price-check column.
$ws.Range("B1").Value = "2025-12-22 04:10"

# This check found the same prices as the previous one (now in column C),
# so copy that snapshot's values/formatting straight across.
$ws.Range("C2:C26").Copy($ws.Range("B2:B26"))
